# ----------------------------------------------------------------------------
# "Setting timeline and update Read me"
#
# 1. Add a new "Timeline" worksheet (as the last tab) with a project agenda
#    table (Front End / Back End / All aspects, due dates, days-left formulas).
# 2. Make the new Timeline sheet the active tab/selection.
# 3. Tweak selections on the CRUD sheet (no more frozen/scrolled topLeftCell,
#    single-cell selection instead of a whole-column selection).
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Create the new "Timeline" sheet as the last tab -----------------------
$newsheet = $wb.Worksheets.Add()
$newsheet.Name = "Timeline"
# Re-fetch the Pages reference AFTER inserting the new sheet (indices shifted)
# then move Timeline after it so Timeline ends up as the last tab.
$pages = $wb.Worksheets.Item("Pages")
$newsheet.Move($null, $pages)
# Sheet references can go stale across structural ops (Add/Move) in this
# engine, so re-fetch a fresh handle by name before using it further.
$timeline = $wb.Worksheets.Item("Timeline")

# --- Header / "days left" reference cell (G1) -------------------------------
$timeline.Range("G1").Formula = "=TODAY()"

# --- Table header row (row 2) ----------------------------------------------
$timeline.Range("A2").Value = "No"
$timeline.Range("B2").Value = "Aspect"
$timeline.Range("C2").Value = "Agenda"
$timeline.Range("D2").Value = "Due Date"
$timeline.Range("E2").Value = "Days Left"

$hdr = $timeline.Range("A2:E2")
$hdr.Font.Name = "Aptos Narrow"
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108   # xlCenter
$hdr.VerticalAlignment = -4108     # xlCenter

# --- Data rows (3-12) --------------------------------------------------------
$timeline.Range("A3").Value = 1
$timeline.Range("B3").Value = "Front End"
$timeline.Range("C3").Value = "Home Page"
$timeline.Range("D3").Value = 45403

$timeline.Range("A4").Value = 2
$timeline.Range("B4").Value = "Front End"
$timeline.Range("C4").Value = "Login / Register Page"
$timeline.Range("D4").Value = 45403

$timeline.Range("A5").Value = 3
$timeline.Range("B5").Value = "Front End"
$timeline.Range("C5").Value = "Product Detail Page"
$timeline.Range("D5").Formula = "=D4+4"

$timeline.Range("A6").Value = 4
$timeline.Range("B6").Value = "Front End"
$timeline.Range("C6").Value = "Admin Page (Change Banner, Product List (Add, Delete), Order List, Change Password)"
$timeline.Range("D6").Formula = "=D5+10"

$timeline.Range("A7").Value = 5
$timeline.Range("B7").Value = "Front End"
$timeline.Range("C7").Value = "Customer Page (My Cart, My Orders)"
$timeline.Range("D7").Formula = "=D6+5"

$timeline.Range("A8").Value = 6
$timeline.Range("B8").Value = "Back End"
$timeline.Range("C8").Value = "Login and Register Process"
$timeline.Range("D8").Formula = "=D7+10"

$timeline.Range("A9").Value = 7
$timeline.Range("B9").Value = "Back End"
$timeline.Range("C9").Value = "Change Password Process"
$timeline.Range("D9").Formula = "=D8+2"

$timeline.Range("A10").Value = 8
$timeline.Range("B10").Value = "Back End"
$timeline.Range("C10").Value = "Product List Process"
$timeline.Range("D10").Formula = "=D9+10"

$timeline.Range("A11").Value = 9
$timeline.Range("B11").Value = "Back End"
$timeline.Range("C11").Value = "Check Out Process"
$timeline.Range("D11").Formula = "=D10+10"

$timeline.Range("A12").Value = 12
$timeline.Range("B12").Value = "All"
$timeline.Range("C12").Value = "Gold Challenge Submission Due Date"
$timeline.Range("D12").Value = 45446

# "Days left" formulas, column E, rows 3-12
for ($row = 3; $row -le 12; $row++) {
    $timeline.Range("E$row").Formula = "=D$row-`$G`$1"
}

# --- Formatting --------------------------------------------------------------
# Whole table: thin border all around, vertically centered text.
$table = $timeline.Range("A2:E12")
$table.Borders.LineStyle = 1
$table.VerticalAlignment = -4108   # xlCenter

# Column B ("Aspect") values are centered.
$timeline.Range("B3:B12").HorizontalAlignment = -4108   # xlCenter

# Column D ("Due Date") gets the custom date format.
$timeline.Range("D3:D12").NumberFormat = "ddd, dd-mmm-yy"

# Re-apply header bold/center/font (kept bold Aptos Narrow, centered).
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4108

# Wrap the long Admin Page agenda text and give that row extra height.
$timeline.Range("C6").WrapText = $true
$timeline.Rows.Item(6).RowHeight = 30

# Distinct left/right-only border treatment on the "Change Password Process"
# cell (matches a stray formatting artifact present in the source workbook).
$timeline.Range("C9").Borders.Item(8).LineStyle = 0    # xlEdgeTop -> none
$timeline.Range("C9").Borders.Item(9).LineStyle = 0    # xlEdgeBottom -> none
$timeline.Range("C9").Borders.Item(7).LineStyle = 1    # xlEdgeLeft -> thin
$timeline.Range("C9").Borders.Item(10).LineStyle = 1   # xlEdgeRight -> thin

# --- Column sizing -----------------------------------------------------------
$timeline.Columns.Item(1).AutoFit()
$timeline.Columns.Item(2).AutoFit()
$timeline.Columns.Item(3).AutoFit()
$timeline.Columns.Item(4).AutoFit()
$timeline.Columns.Item(5).ColumnWidth = 11.85546875
$timeline.Columns.Item(6).ColumnWidth = 10.42578125
$timeline.Columns.Item(7).ColumnWidth = 13.140625

# --- Selections ----------------------------------------------------------
$crud = $wb.Worksheets.Item("CRUD")
$crud.Activate()
$crud.Range("E13").Select()

# Make Timeline the active tab last, so its selection + zoom stick as the
# saved view (this engine applies ActiveWindow.Zoom to whichever sheet is
# currently active, so Zoom must be set AFTER activating Timeline).
$timeline = $wb.Worksheets.Item("Timeline")
$timeline.Activate()
$timeline.Range("G2").Select()
$timeline.Application.ActiveWindow.Zoom = 145

Write-Host "Timeline sheet created."
